$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-04-07 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-04-08 Tuesday", 2)

# New arithmetic-problem values for the 20x5 practice table, listed in
# row-major order (row 1 col 1, row 1 col 2, ... row 1 col 5, row 2 col 1, ...).
$newValues = @(
    "80-26=",
    "40+27=",
    "24+60=",
    "31+53=",
    "92-47=",
    "56+19=",
    "36+51=",
    "63-53=",
    "67+2=",
    "19+33=",
    "0+88=",
    "12+85=",
    "10+64=",
    "57-55=",
    "72-58=",
    "45-19=",
    "58-57=",
    "1+8=",
    "78-23=",
    "92-84=",
    "41+27=",
    "66-1=",
    "63+34=",
    "30+44=",
    "66-26=",
    "29+10=",
    "11+4=",
    "10+53=",
    "60+9=",
    "64+15=",
    "96-36=",
    "75-62=",
    "14+59=",
    "24+20=",
    "21+49=",
    "48+31=",
    "80-14=",
    "12+51=",
    "36+0=",
    "23+0=",
    "29+18=",
    "7+33=",
    "81-11=",
    "99-57=",
    "55+35=",
    "67+32=",
    "8+67=",
    "64-31=",
    "30-15=",
    "29+29=",
    "5+94=",
    "5+43=",
    "44+47=",
    "58-34=",
    "69-53=",
    "82-1=",
    "31-27=",
    "78-35=",
    "25+12=",
    "55-44=",
    "86-33=",
    "54-41=",
    "42+33=",
    "30-5=",
    "17+59=",
    "68-56=",
    "77-4=",
    "95-77=",
    "49-2=",
    "86-26=",
    "64-21=",
    "55+10=",
    "8+35=",
    "32+25=",
    "58+33=",
    "50+27=",
    "5+13=",
    "53+11=",
    "8+27=",
    "85-23=",
    "10+75=",
    "37+6=",
    "58-16=",
    "28+8=",
    "49+3=",
    "86-83=",
    "88-51=",
    "76-30=",
    "20+74=",
    "2+96=",
    "39-33=",
    "88-17=",
    "41+1=",
    "98-29=",
    "56-48=",
    "97-37=",
    "96-59=",
    "19+76=",
    "31+60=",
    "63+26="
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}

Write-Output ("Updated " + $i.ToString() + " cells")
